# BOM.xlsx update - "Components added to BOM"
#
# Adds/updates several component rows in the Bill Of Materials sheet:
#  - row 5  (L298 DC Motor driver): quantity 1 -> 2, adds hyperlink on D5
#  - row 10 (breadboard): component text + link replaced, adds hyperlink on D10
#  - row 12 (LiPo battery): component text + link replaced, adds hyperlink on D12
#  - rows 14-18: five new component rows appended (row 15 gets a hyperlink on D15)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 5: L298 DC Motor driver - quantity bump + new hyperlink on D5 ----
$ws.Range("C5").Value = 2

$ws.Hyperlinks.Add($ws.Range("D5"), "https://goo.gl/TzPNkm")
$ws.Range("D5").Style = "Hyperlink"

# ---- Row 10: breadboard description/link replaced ----
$ws.Range("B10").Value = "Bread board with Power module and jumper cables"
$ws.Range("D10").Value = "https://goo.gl/GiHjuM"
$ws.Hyperlinks.Add($ws.Range("D10"), "https://goo.gl/GiHjuM")
$ws.Range("D10").Style = "Hyperlink"

# ---- Row 12: LiPo battery description/link replaced ----
$ws.Range("B12").Value = "7.4V 1500 mAh 25C LiPo battery for motors"
$ws.Range("D12").Value = "https://goo.gl/PEDWJU"
$ws.Hyperlinks.Add($ws.Range("D12"), "https://goo.gl/PEDWJU")
$ws.Range("D12").Style = "Hyperlink"

# ---- New rows 14-18 ----

# Row 14: 5v 2.4 A Power bank for Raspberry Pi (plain url text, like row 6)
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "5v 2.4 A Power bank for Raspberry Pi"
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = "https://goo.gl/o323cW"
$ws.Range("C6").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D14").PasteSpecial(-4122)

# Row 15: 9V 1A AC-DC Converter for Breadboard (hyperlinked)
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "9V 1A AC-DC Converter for Breadboard"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = "https://goo.gl/Jt5Dpn"
$ws.Range("C6").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("D15"), "https://goo.gl/Jt5Dpn")
$ws.Range("D15").Style = "Hyperlink"

# Row 16: Double sided foam tape
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Double sided foam tape"
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = "https://goo.gl/7ENrhd"
$ws.Range("C6").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D16").PasteSpecial(-4122)

# Row 17: Anti Static ESD Gloves
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "Anti Static ESD Gloves"
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = "https://goo.gl/7yG2hZ"
$ws.Range("C6").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D17").PasteSpecial(-4122)

# Row 18: Raspberry Pi base plate Holder
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "Raspberry Pi base plate Holder"
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = "https://goo.gl/eMxNbV"
$ws.Range("C6").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D18").PasteSpecial(-4122)

# Restore the active selection to where the author ended up
$ws.Range("B22").Select()
